$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

$newValues = @("CIGNA", "CS2FEUR", "CS2FUSD", "DL4COINV", "DL4COINVEUR", "SMA ILM", "SMA SUMI")

$row = 55
foreach ($val in $newValues) {
    $ws2.Cells.Item($row, 1).Value = $val
    $row = $row + 1
}

# Copy the formatting from the last pre-existing row so the new rows
# pick up the same cell style (fill) as the rest of column A.
$src = $ws2.Range("A54")
$dst = $ws2.Range("A55:A61")
$src.Copy()
$dst.PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Reflect the new list bounds in the view (scroll/selection) similar to
# what Excel records after editing near the bottom of the list, then
# restore Sheet1 as the active/visible tab (Sheet2 is veryHidden).
$ws2.Range("A55:A61").Select() | Out-Null
$ws1.Select() | Out-Null
